$wb = $excel.ActiveWorkbook

# --- Rename header in "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename header in "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header style (bold, centered, bordered) from "Weekly Quantity"!A1:B1
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the existing date-format style from "Weekly Quantity"!A2 for column A of the new sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A34").PasteSpecial(-4122)

$data = @(
    @(44955.99999999999, 11, 0.02334687334229298, 21.14339148652448),
    @(44969.99999999999, 11, 0.8506102939307333, 21.57115065992996),
    @(44976.99999999999, 11, 1.228005154226867, 21.47387447248423),
    @(44983.99999999999, 11, 1.575907928045903, 22.13294805180795),
    @(44990.99999999999, 11, 1.633744029829142, 21.32364241973748),
    @(44997.99999999999, 12, 0.3774811764036823, 21.87832089809059),
    @(45004.99999999999, 12, 1.996007286820381, 22.24384166000689),
    @(45011.99999999999, 12, 0.8673032855542968, 22.39083455644161),
    @(45018.99999999999, 12, 1.644801668159787, 22.35158111629039),
    @(45039.99999999999, 12, 2.293646032627447, 22.1889627684166),
    @(45067.99999999999, 13, 2.869672526916099, 22.29930586663795),
    @(45074.99999999999, 13, 2.56359993027601, 22.89169399907818),
    @(45081.99999999999, 13, 3.185053547520376, 22.32120770124208),
    @(45088.99999999999, 13, 2.895520675787333, 23.29415696614171),
    @(45102.99999999999, 13, 1.9028486899069, 23.2112499413412),
    @(45109.99999999999, 13, 3.261560495392179, 24.11133054583862),
    @(45116.99999999999, 13, 2.816946863620842, 23.0077209854293),
    @(45123.99999999999, 13, 4.089647469387242, 23.34940874707245),
    @(45130.99999999999, 13, 2.804922551601836, 22.63479169913035),
    @(45165.99999999999, 14, 3.58696414867147, 23.25726964897562),
    @(45172.99999999999, 14, 3.813776293098318, 24.5973668521694),
    @(45179.99999999999, 14, 4.207495131605745, 24.51966407072264),
    @(45186.99999999999, 14, 4.353944059862345, 24.3421583033907),
    @(45193.99999999999, 14, 4.170992108519796, 24.74805177504551),
    @(45200.99999999999, 14, 4.644400943768333, 24.77093058156694),
    @(45207.99999999999, 15, 3.814516820498381, 24.82080281279284),
    @(45214.99999999999, 15, 4.381618715593109, 25.0512717297256),
    @(45221.99999999999, 15, 4.760034165749527, 23.96401285890631),
    @(45228.99999999999, 15, 4.854169668957112, 24.58242465072745),
    @(45235.99999999999, 15, 4.385029413815691, 25.16360358767011),
    @(45242.99999999999, 15, 5.095640454923153, 25.69291376582676),
    @(45249.99999999999, 15, 4.323423001167821, 26.2799465896371),
    @(45256.99999999999, 15, 4.775076802880976, 25.21724367851822)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $wsForecast.Cells.Item($row, 1).Value = $vals[0]
    $wsForecast.Cells.Item($row, 2).Value = $vals[1]
    $wsForecast.Cells.Item($row, 3).Value = $vals[2]
    $wsForecast.Cells.Item($row, 4).Value = $vals[3]
}

Write-Host "done"
